# Daily attendance processing - 2025-10-31 07:42:46
# Reorders the "Recorded By" (column G) contributor lists so the
# backup/admin addresses are listed first, matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "backup@backdoor.com, system, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    7   = "System, admin@admin.com"
    8   = "backup@backdoor.com, System"
    11  = "dnasr281@gmail.com, System"
    17  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, system, System"
    31  = "backup@backdoor.com, System"
    32  = "backup@backdoor.com, System"
    34  = "System, admin@admin.com"
    35  = "backup@backdoor.com, System"
    38  = "dnasr281@gmail.com, System"
    44  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, system, System"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    61  = "System, admin@admin.com"
    62  = "backup@backdoor.com, System"
    65  = "dnasr281@gmail.com, System"
    71  = "dnasr281@gmail.com, System"
    83  = "backup@backdoor.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    90  = "dnasr281@gmail.com, admin@admin.com"
    96  = "dnasr281@gmail.com, System"
    97  = "dnasr281@gmail.com, System"
    99  = "dnasr281@gmail.com, System"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    116 = "dnasr281@gmail.com, admin@admin.com"
    122 = "dnasr281@gmail.com, System"
    123 = "dnasr281@gmail.com, System"
    125 = "dnasr281@gmail.com, System"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    142 = "dnasr281@gmail.com, admin@admin.com"
    148 = "dnasr281@gmail.com, System"
    149 = "dnasr281@gmail.com, System"
    151 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
